# Apply targeted cell updates to match the committed diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 40
$ws.Range("C9").Value = 2
$ws.Range("G9").Value = "'512.00"
$ws.Range("C10").Value = 1
$ws.Range("G10").Value = "'472.00"
$ws.Range("C11").Value = 76
$ws.Range("G11").Value = "'50312.00"
$ws.Range("A12").Value = "'P. point"
$ws.Range("C12").Value = 95
$ws.Range("D12").Value = "'6"
$ws.Range("E12").Value = "'On board"
$ws.Range("F12").Value = 136
$ws.Range("G12").Value = "'12920.00"
$ws.Range("C13").Value = 82
$ws.Range("G13").Value = "'1886.00"
$ws.Range("C14").Value = 42
$ws.Range("G14").Value = "'2100.00"
$ws.Range("C15").Value = 41
$ws.Range("D15").Value = "'6.0"
$ws.Range("E15").Value = "'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F15").Value = 78
$ws.Range("G15").Value = "'3198.00"
$ws.Range("C16").Value = 76
$ws.Range("D16").Value = "'8.0"
$ws.Range("E16").Value = "'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F16").Value = 30
$ws.Range("G16").Value = "'2280.00"
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = "'9.0"
$ws.Range("E17").Value = "'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F17").Value = 219
$ws.Range("G17").Value = "'1095.00"
$ws.Range("C18").Value = 24
$ws.Range("D18").Value = "'10.0"
$ws.Range("E18").Value = "'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F18").Value = 303
$ws.Range("G18").Value = "'7272.00"
$ws.Range("C19").Value = 73
$ws.Range("C20").Value = 64
$ws.Range("G20").Value = "'2560.00"
$ws.Range("C21").Value = 44
$ws.Range("G21").Value = "'2464.00"
$ws.Range("C22").Value = 71
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = "'19"
$ws.Range("E23").Value = "'2 x 2.5 sq. mm. + 1x1.5sqmm"
$ws.Range("F23").Value = 81
$ws.Range("G23").Value = "'162.00"
$ws.Range("A24").Value = "'Mtr."
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = "'20"
$ws.Range("E24").Value = "'2 x 4.0 sq. mm. + 1 x 2.5 sq. mm."
$ws.Range("F24").Value = 122
$ws.Range("G24").Value = "'2806.00"
$ws.Range("A25").Value = "'Set"
$ws.Range("C25").Value = 37
$ws.Range("D25").Value = "'13.0"
$ws.Range("E25").Value = "'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("F25").Value = 5733
$ws.Range("G25").Value = "'212121.00"
$ws.Range("A26").Value = "'"
$ws.Range("C26").Value = 68
$ws.Range("D26").Value = "'15.0"
$ws.Range("E26").Value = "'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = "'0.00"
$ws.Range("A27").Value = "'"
$ws.Range("C27").Value = 31
$ws.Range("D27").Value = "'16.0"
$ws.Range("E27").Value = "'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = "'0.00"
$ws.Range("A28").Value = "'Each"
$ws.Range("C28").Value = 98
$ws.Range("D28").Value = "'27"
$ws.Range("E28").Value = "'1170mm(+/-10%) LED batten with min. lumen output 2200 lm"
$ws.Range("F28").Value = 492
$ws.Range("G28").Value = "'48216.00"
$ws.Range("A29").Value = "'"
$ws.Range("C29").Value = 48
$ws.Range("D29").Value = "'17.0"
$ws.Range("E29").Value = "'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = "'0.00"
$ws.Range("C30").Value = 93
$ws.Range("D30").Value = "'29"
$ws.Range("E30").Value = "'Single pole MCB   (With B/C curve tripping Characteristics)"
$ws.Range("A31").Value = "'Each"
$ws.Range("C31").Value = 41
$ws.Range("D31").Value = "'30"
$ws.Range("E31").Value = "' 6 A to 32 A rating"
$ws.Range("F31").Value = 187
$ws.Range("G31").Value = "'7667.00"
$ws.Range("C32").Value = 93
$ws.Range("C33").Value = 77
$ws.Range("G33").Value = "'69300.00"
$ws.Range("C34").Value = 60
$ws.Range("C35").Value = 68
$ws.Range("B36").Value = 0
$ws.Range("C36").Value = 34
$ws.Range("D36").Value = "'36"
$ws.Range("E36").Value = "'Total"
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = "'0.00"
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = "'"
$ws.Range("A37").Value = "'%"
$ws.Range("B37").Value = 0
$ws.Range("C37").Value = 22
$ws.Range("D37").Value = "'37"
$ws.Range("E37").Value = "'Add Tender Premium "
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = "'0.00"
$ws.Range("H37").Value = 0
$ws.Range("E39").Value = "'Grand Total Rs."
$ws.Range("G39").Value = "'427343.00"
$ws.Range("H39").Value = "'427343.00"
$ws.Range("A40").Value = "'"
$ws.Range("B40").Value = "'"
$ws.Range("C40").Value = "'"
$ws.Range("D40").Value = "'"
$ws.Range("E40").Value = "'Tender Premium @ 0%"
$ws.Range("F40").Value = "'"
$ws.Range("G40").Value = "'0.00"
$ws.Range("H40").Value = "'0.00"
$ws.Range("I40").Value = "'"
$ws.Range("A41").Value = "'"
$ws.Range("B41").Value = "'"
$ws.Range("C41").Value = "'"
$ws.Range("D41").Value = "'"
$ws.Range("E41").Value = "'NET PAYABLE AMOUNT Rs."
$ws.Range("F41").Value = "'"
$ws.Range("G41").Value = "'427343.00"
$ws.Range("H41").Value = "'427343.00"
$ws.Range("I41").Value = "'"

# Row 38 loses its former "Tender Premium" content (cells B38:I38 are removed entirely).
$ws.Range("B38:I38").ClearContents()

